$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 1.593927
$ws.Range("H2").Value = 3.187854
$ws.Range("I2").Value = 0.1926821744909273
$ws.Range("J2").Value = 0.1419711187769723
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.889399666666667
$ws.Range("N2").Value = 17.668199
$ws.Range("O2").Value = 0.8160192454225522
$ws.Range("P2").Value = 0.8160192454225521
$ws.Range("Q2").Value = 9.387273142490999
$ws.Range("R2").Value = 56.323638854946
$ws.Range("S2").Value = 0.157232362634463
$ws.Range("T2").Value = 0.1158511652161805

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1.593927
$ws.Range("H3").Value = 3.187854
$ws.Range("I3").Value = 0.1926821744909273
$ws.Range("J3").Value = 0.1419711187769723
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.327831666666667
$ws.Range("N3").Value = 3.983495
$ws.Range("O3").Value = 0.1839807545774479
$ws.Range("P3").Value = 0.1839807545774478
$ws.Range("Q3").Value = 2.116466744955
$ws.Range("R3").Value = 12.69880046973
$ws.Range("S3").Value = 0.03544981185646428
$ws.Range("T3").Value = 0.02611995356079184

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.813729
$ws.Range("H4").Value = 2.441187
$ws.Range("I4").Value = 0.09836778796414627
$ws.Range("J4").Value = 0.1087182943553252
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.889399666666667
$ws.Range("N4").Value = 17.668199
$ws.Range("O4").Value = 0.8160192454225522
$ws.Range("P4").Value = 0.8160192454225521
$ws.Range("Q4").Value = 4.792375301357001
$ws.Range("R4").Value = 43.13137771221301
$ws.Range("S4").Value = 0.08027000810838826
$ws.Range("T4").Value = 0.08871622052345936

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.813729
$ws.Range("H5").Value = 2.441187
$ws.Range("I5").Value = 0.09836778796414627
$ws.Range("J5").Value = 0.1087182943553252
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.327831666666667
$ws.Range("N5").Value = 3.983495
$ws.Range("O5").Value = 0.1839807545774479
$ws.Range("P5").Value = 0.1839807545774478
$ws.Range("Q5").Value = 1.080495134285
$ws.Range("R5").Value = 9.724456208565002
$ws.Range("S5").Value = 0.01809777985575802
$ws.Range("T5").Value = 0.02000207383186581

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.145019
$ws.Range("H6").Value = 6.435057
$ws.Range("I6").Value = 0.2593010377792423
$ws.Range("J6").Value = 0.2865853460301467
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.889399666666667
$ws.Range("N6").Value = 17.668199
$ws.Range("O6").Value = 0.8160192454225522
$ws.Range("P6").Value = 0.8160192454225521
$ws.Range("Q6").Value = 12.63287418359367
$ws.Range("R6").Value = 113.695867652343
$ws.Range("S6").Value = 0.211594637185902
$ws.Range("T6").Value = 0.2338591578166813

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.145019
$ws.Range("H7").Value = 6.435057
$ws.Range("I7").Value = 0.2593010377792423
$ws.Range("J7").Value = 0.2865853460301467
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.327831666666667
$ws.Range("N7").Value = 3.983495
$ws.Range("O7").Value = 0.1839807545774479
$ws.Range("P7").Value = 0.1839807545774478
$ws.Range("Q7").Value = 2.848224153801667
$ws.Range("R7").Value = 25.63401738421501
$ws.Range("S7").Value = 0.04770640059334032
$ws.Range("T7").Value = 0.05272618821346539

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.113821666666667
$ws.Range("H8").Value = 3.341465
$ws.Range("I8").Value = 0.1346445481684181
$ws.Range("J8").Value = 0.1488121866321657
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.889399666666667
$ws.Range("N8").Value = 17.668199
$ws.Range("O8").Value = 0.8160192454225522
$ws.Range("P8").Value = 0.8160192454225521
$ws.Range("Q8").Value = 6.559740952392779
$ws.Range("R8").Value = 59.03766857153501
$ws.Range("S8").Value = 0.109872542596653
$ws.Range("T8").Value = 0.1214336082452598

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.113821666666667
$ws.Range("H9").Value = 3.341465
$ws.Range("I9").Value = 0.1346445481684181
$ws.Range("J9").Value = 0.1488121866321657
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.327831666666667
$ws.Range("N9").Value = 3.983495
$ws.Range("O9").Value = 0.1839807545774479
$ws.Range("P9").Value = 0.1839807545774478
$ws.Range("Q9").Value = 1.478967680019445
$ws.Range("R9").Value = 13.310709120175
$ws.Range("S9").Value = 0.02477200557176509
$ws.Range("T9").Value = 0.02737857838690584

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.83705
$ws.Range("H10").Value = 5.51115
$ws.Range("I10").Value = 0.222072145492584
$ws.Range("J10").Value = 0.2454391359352439
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.889399666666667
$ws.Range("N10").Value = 17.668199
$ws.Range("O10").Value = 0.8160192454225522
$ws.Range("P10").Value = 0.8160192454225521
$ws.Range("Q10").Value = 10.81912165765
$ws.Range("R10").Value = 97.37209491885001
$ws.Range("S10").Value = 0.1812151445942256
$ws.Range("T10").Value = 0.2002830585030409

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.83705
$ws.Range("H11").Value = 5.51115
$ws.Range("I11").Value = 0.222072145492584
$ws.Range("J11").Value = 0.2454391359352439
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.327831666666667
$ws.Range("N11").Value = 3.983495
$ws.Range("O11").Value = 0.1839807545774479
$ws.Range("P11").Value = 0.1839807545774478
$ws.Range("Q11").Value = 2.43929316325
$ws.Range("R11").Value = 21.95363846925
$ws.Range("S11").Value = 0.0408570008983584
$ws.Range("T11").Value = 0.04515607743220297

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.768765
$ws.Range("H12").Value = 1.53753
$ws.Range("I12").Value = 0.09293230610468217
$ws.Range("J12").Value = 0.0684739182701461
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.889399666666667
$ws.Range("N12").Value = 17.668199
$ws.Range("O12").Value = 0.8160192454225522
$ws.Range("P12").Value = 0.8160192454225521
$ws.Range("Q12").Value = 4.527564334745001
$ws.Range("R12").Value = 27.16538600847
$ws.Range("S12").Value = 0.07583455030292038
$ws.Range("T12").Value = 0.05587603511793012

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.768765
$ws.Range("H13").Value = 1.53753
$ws.Range("I13").Value = 0.09293230610468217
$ws.Range("J13").Value = 0.0684739182701461
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.327831666666667
$ws.Range("N13").Value = 3.983495
$ws.Range("O13").Value = 0.1839807545774479
$ws.Range("P13").Value = 0.1839807545774478
$ws.Range("Q13").Value = 1.020790511225
$ws.Range("R13").Value = 6.124743067350001
$ws.Range("S13").Value = 0.01709775580176179
$ws.Range("T13").Value = 0.01259788315221597
